# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff"
# - Latest Handoff/Handback style datetime stamps are refreshed
# - The "Status" columns are widened slightly (to fit the longer text)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps ---
$wsOverview.Range("G2").Value = "2016-09-01 15:23:53"
$wsZhCn.Range("H2").Value     = "2016-09-01 15:23:47"
$wsDeDe.Range("H2").Value     = "2016-09-01 15:23:53"

# --- Widen the Status columns to fit the new text ---
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3
$wsZhCn.Range("C1").ColumnWidth     = 16.3
$wsDeDe.Range("C1").ColumnWidth     = 16.3
